# Update "想去人数" (column F) figures across the four worksheets to match
# the refreshed data snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 1269
    5  = 5638
    6  = 1803
    7  = 1803
    8  = 6372
    9  = 140
    10 = 1938
    11 = 519
    12 = 14
    14 = 34
    18 = 7971
    19 = 7971
    22 = 185
    24 = 1748
    25 = 848
    26 = 4
    31 = 1780
    32 = 807
    33 = 380
    36 = 6
    38 = 88
    39 = 3922
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    4 = 368
    8 = 1
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2Updates[$row]
}

# --- Sheet "本地生活" (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$sheet3Updates = @{
    2 = 9543
    4 = 689
    5 = 277
}
foreach ($row in $sheet3Updates.Keys) {
    $ws3.Cells.Item($row, 6).Value = $sheet3Updates[$row]
}

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 9543
    4  = 689
    5  = 1269
    9  = 368
    10 = 5638
    11 = 277
    12 = 1803
    13 = 1804
    14 = 6372
    15 = 140
    16 = 1938
    18 = 519
    21 = 34
    24 = 7971
    25 = 7971
    28 = 185
    30 = 1748
    31 = 848
    32 = 4
    36 = 1780
    37 = 807
    39 = 380
    46 = 88
    47 = 3922
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
